$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, bordered, centered) from H1 into the two new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for I2:J26
$data = @(
    @(2, 9, 9),
    @(3, 6, 6),
    @(4, 10, 11),
    @(5, 8, 8),
    @(6, 7, 8),
    @(7, 8, 8),
    @(8, 9, 9),
    @(9, 9, 9),
    @(10, 8, 8),
    @(11, 6, 6),
    @(12, 7, 7),
    @(13, 5, 5),
    @(14, 8, 8),
    @(15, 7, 7),
    @(16, 6, 6),
    @(17, 7, 7),
    @(18, 6, 6),
    @(19, 5, 5),
    @(20, 9, 9),
    @(21, 9, 9),
    @(22, 10, 10),
    @(23, 7, 7),
    @(24, 7, 7),
    @(25, 4, 4),
    @(26, 4, 4)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}

